# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 12918
$ws1.Range("F4").Value  = 28
$ws1.Range("F6").Value  = 76
$ws1.Range("F8").Value  = 21
$ws1.Range("F9").Value  = 17
$ws1.Range("F10").Value = 12865
$ws1.Range("F11").Value = 285
$ws1.Range("F12").Value = 39
$ws1.Range("F13").Value = 8691
$ws1.Range("F14").Value = 7689
$ws1.Range("F15").Value = 197
$ws1.Range("F16").Value = 112
$ws1.Range("F18").Value = 126
$ws1.Range("F21").Value = 15

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 12918
$ws4.Range("F5").Value  = 28
$ws4.Range("F7").Value  = 76
$ws4.Range("F9").Value  = 21
$ws4.Range("F10").Value = 17
$ws4.Range("F11").Value = 12865
$ws4.Range("F12").Value = 285
$ws4.Range("F13").Value = 39
$ws4.Range("F14").Value = 8691
$ws4.Range("F15").Value = 7689
$ws4.Range("F16").Value = 197
$ws4.Range("F17").Value = 112
$ws4.Range("F19").Value = 126
$ws4.Range("F22").Value = 15

$wb.Save()
